# Generate Report for Handback
# Updates the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) timestamps on the
# per-language report sheets, as produced by a fresh report run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-20 14:23:35"
$wsZhCn.Range("H2").Value = "2016-03-20 14:23:56"
$wsZhCn.Range("E4").Value = "2016-03-20 14:23:35"
$wsZhCn.Range("H4").Value = "2016-03-20 14:23:56"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-20 14:23:38"
$wsDeDe.Range("H2").Value = "2016-03-20 14:24:02"
$wsDeDe.Range("E4").Value = "2016-03-20 14:23:38"
$wsDeDe.Range("H4").Value = "2016-03-20 14:24:02"
